$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.124482532946104
$ws.Range("D2").Value = 0.006277004901020788
$ws.Range("E2").Value = 0.6738359881879035
$ws.Range("F2").Value = 0.9506464327876358
$ws.Range("G2").Value = 0.8596174217706221
$ws.Range("H2").Value = 0.7517988108240559
$ws.Range("L2").Value = 0.2489952740098573
$ws.Range("M2").Value = 0.2683719018378241
$ws.Range("N2").Value = 1.745962422070548
$ws.Range("B3").Value = 1.080744149997571
$ws.Range("D3").Value = 0.006315920974921951
$ws.Range("E3").Value = 0.6412361146004315
$ws.Range("F3").Value = 0.8973036700973154
$ws.Range("G3").Value = 0.7987637060012958
$ws.Range("H3").Value = 0.7295980631369048
$ws.Range("L3").Value = 0.2240916001475739
$ws.Range("M3").Value = 0.2519584832880852
$ws.Range("N3").Value = 1.752865445488439
$ws.Range("B4").Value = 1.05437849822485
$ws.Range("D4").Value = 0.006349443960331769
$ws.Range("E4").Value = 0.6210873244879167
$ws.Range("F4").Value = 0.8653051507813103
$ws.Range("G4").Value = 0.7620875950661627
$ws.Range("H4").Value = 0.7165442534656563
$ws.Range("L4").Value = 0.2088521690103704
$ws.Range("M4").Value = 0.241983852367774
$ws.Range("N4").Value = 1.757842727132115
$ws.Range("B5").Value = 1.043757879920179
$ws.Range("D5").Value = 0.006365491613628649
$ws.Range("E5").Value = 0.6128432952560274
$ws.Range("F5").Value = 0.8524530781679118
$ws.Range("G5").Value = 0.7473124001449492
$ws.Range("H5").Value = 0.7113689475358171
$ws.Range("L5").Value = 0.202654940261155
$ws.Range("M5").Value = 0.2379451547681128
$ws.Range("N5").Value = 1.760057309487351
$ws.Range("B6").Value = 1.042001811694092
$ws.Range("D6").Value = 0.006368299469219352
$ws.Range("E6").Value = 0.6114723747047535
$ws.Range("F6").Value = 0.8503302711721688
$ws.Range("G6").Value = 0.7448692234748648
$ws.Range("H6").Value = 0.7105182738782787
$ws.Range("L6").Value = 0.2016266784135183
$ws.Range("M6").Value = 0.2372761075080234
$ws.Range("N6").Value = 1.760436306786602
$ws.Range("B7").Value = 1.054234763862667
$ws.Range("D7").Value = 0.006349650766132697
$ws.Range("E7").Value = 0.6209762771125824
$ws.Range("F7").Value = 0.8651310664297682
$ws.Range("G7").Value = 0.7618876440430711
$ws.Range("H7").Value = 0.7164738748381012
$ws.Range("L7").Value = 0.2087685385167788
$ws.Range("M7").Value = 0.2419292795210382
$ws.Range("N7").Value = 1.757871838751129
$ws.Range("B8").Value = 1.109300090248041
$ws.Range("D8").Value = 0.006288404024370919
$ws.Range("E8").Value = 0.6626230328141389
$ws.Range("F8").Value = 0.9320960654694517
$ws.Range("G8").Value = 0.838490693219029
$ws.Range("H8").Value = 0.7440235119812257
$ws.Range("L8").Value = 0.2403977881022854
$ws.Range("M8").Value = 0.2626911599767467
$ws.Range("N8").Value = 1.74818955595147
$ws.Range("B9").Value = 1.221158952951413
$ws.Range("D9").Value = 0.006246258354522638
$ws.Range("E9").Value = 0.7432443766454213
$ws.Range("F9").Value = 1.069498704424362
$ws.Range("G9").Value = 0.9942902112834133
$ws.Range("H9").Value = 0.802679884455074
$ws.Range("L9").Value = 0.302834661549042
$ws.Range("M9").Value = 0.3042236540089647
$ws.Range("N9").Value = 1.735043010007786
$ws.Range("B10").Value = 1.305699800385696
$ws.Range("D10").Value = 0.006265000785120378
$ws.Range("E10").Value = 0.8018471843548411
$ws.Range("F10").Value = 1.174308494429624
$ws.Range("G10").Value = 1.112336544236427
$ws.Range("H10").Value = 0.8486720802145555
$ws.Range("L10").Value = 0.3489676429383621
$ws.Range("M10").Value = 0.335238912805643
$ws.Range("N10").Value = 1.728917556302392
$ws.Range("B11").Value = 1.344671434278837
$ws.Range("D11").Value = 0.006284797858114644
$ws.Range("E11").Value = 0.8283728237049672
$ws.Range("F11").Value = 1.222860840741504
$ws.Range("G11").Value = 1.166855732849285
$ws.Range("H11").Value = 0.8702406906226372
$ws.Range("L11").Value = 0.3700137260755127
$ws.Range("M11").Value = 0.3494581886442774
$ws.Range("N11").Value = 1.72689270122909
$ws.Range("B12").Value = 1.359502600830353
$ws.Range("D12").Value = 0.006293954788379352
$ws.Range("E12").Value = 0.8383983392071457
$ws.Range("F12").Value = 1.241374524414681
$ws.Range("G12").Value = 1.187621526954899
$ws.Range("H12").Value = 0.8785023730177386
$ws.Range("L12").Value = 0.3779920470697959
$ws.Range("M12").Value = 0.3548585067183936
$ws.Range("N12").Value = 1.726235009866002
$ws.Range("B13").Value = 1.356305183078575
$ws.Range("D13").Value = 0.006291908207820995
$ws.Range("E13").Value = 0.836240018017449
$ws.Range("F13").Value = 1.237381543159472
$ws.Range("G13").Value = 1.183143835688099
$ws.Range("H13").Value = 0.8767188675459749
$ws.Range("L13").Value = 0.3762733891261973
$ws.Range("M13").Value = 0.3536947505692183
$ws.Range("N13").Value = 1.726371811219465
$ws.Range("B14").Value = 1.345890133387115
$ws.Range("D14").Value = 0.006285517702870891
$ws.Range("E14").Value = 0.8291980132187717
$ws.Range("F14").Value = 1.224381394753152
$ws.Range("G14").Value = 1.168561715196432
$ws.Range("H14").Value = 0.8709184913212766
$ws.Range("L14").Value = 0.3706699344243418
$ws.Range("M14").Value = 0.3499021600657528
$ws.Range("N14").Value = 1.726836409012904
$ws.Range("B15").Value = 1.339520169016907
$ws.Range("D15").Value = 0.006281820722616516
$ws.Range("E15").Value = 0.824882084828559
$ws.Range("F15").Value = 1.216435162031217
$ws.Range("G15").Value = 1.159645527810341
$ws.Range("H15").Value = 0.8673778846171558
$ws.Range("L15").Value = 0.3672387802016033
$ws.Range("M15").Value = 0.3475811427368143
$ws.Range("N15").Value = 1.727135180593663
$ws.Range("B16").Value = 1.303163212823108
$ws.Range("D16").Value = 0.006263937162252375
$ws.Range("E16").Value = 0.8001109906841037
$ws.Range("F16").Value = 1.171153264954739
$ws.Range("G16").Value = 1.108790307711956
$ws.Range("H16").Value = 0.8472756283894114
$ws.Range("L16").Value = 0.3475934386626136
$ws.Range("M16").Value = 0.3343118626412647
$ws.Range("N16").Value = 1.729065173193732
$ws.Range("B17").Value = 1.280990652350624
$ws.Range("D17").Value = 0.006255880871684383
$ws.Range("E17").Value = 0.7848806177081684
$ws.Range("F17").Value = 1.143599733322191
$ws.Range("G17").Value = 1.077804067772206
$ws.Range("H17").Value = 0.8351099489231899
$ws.Range("L17").Value = 0.3355570160709931
$ws.Range("M17").Value = 0.3261997826896561
$ws.Range("N17").Value = 1.730443915736004
$ws.Range("B18").Value = 1.268285945714126
$ws.Range("D18").Value = 0.006252305601957886
$ws.Range("E18").Value = 0.7761079860608504
$ws.Range("F18").Value = 1.127833781900449
$ws.Range("G18").Value = 1.060058554608077
$ws.Range("H18").Value = 0.8281733568235836
$ws.Range("L18").Value = 0.3286396271642502
$ws.Range("M18").Value = 0.3215443129094098
$ws.Range("N18").Value = 1.73130864760239
$ws.Range("B19").Value = 1.263992669086235
$ws.Range("D19").Value = 0.006251275680664747
$ws.Range("E19").Value = 0.7731355738753933
$ws.Range("F19").Value = 1.122509734489398
$ws.Range("G19").Value = 1.054063362299189
$ws.Range("H19").Value = 0.8258351510678494
$ws.Range("L19").Value = 0.3262984872898471
$ws.Range("M19").Value = 0.3199698365184034
$ws.Range("N19").Value = 1.731613761940324
$ws.Range("B20").Value = 1.283345956708501
$ws.Range("D20").Value = 0.006256628643928508
$ws.Range("E20").Value = 0.7865032133568235
$ws.Range("F20").Value = 1.14652433898641
$ws.Range("G20").Value = 1.081094621840151
$ws.Range("H20").Value = 0.8363987069951975
$ws.Range("L20").Value = 0.3368377297591962
$ws.Range("M20").Value = 0.3270622532633496
$ws.Range("N20").Value = 1.730289727211769
$ws.Range("B21").Value = 1.348947295536959
$ws.Range("D21").Value = 0.006287349375391216
$ws.Range("E21").Value = 0.8312669411732116
$ws.Range("F21").Value = 1.228196362059975
$ws.Range("G21").Value = 1.172841543473083
$ws.Range("H21").Value = 0.8726196388731466
$ws.Range("L21").Value = 0.3723155707328942
$ws.Range("M21").Value = 0.3510157079297329
$ws.Range("N21").Value = 1.726696988556654
$ws.Range("B22").Value = 1.392249533468032
$ws.Range("D22").Value = 0.006317122779787354
$ws.Range("E22").Value = 0.8604108878595298
$ws.Range("F22").Value = 1.282320710571526
$ws.Range("G22").Value = 1.233507473863654
$ws.Range("H22").Value = 0.8968412451160361
$ws.Range("L22").Value = 0.3955526756291476
$ws.Range("M22").Value = 0.3667626626139224
$ws.Range("N22").Value = 1.724984512419155
$ws.Range("B23").Value = 1.369099284040828
$ws.Range("D23").Value = 0.006300331838033912
$ws.Range("E23").Value = 0.8448664475159688
$ws.Range("F23").Value = 1.25336441154289
$ws.Range("G23").Value = 1.201063593780447
$ws.Range("H23").Value = 0.8838630892093136
$ws.Range("L23").Value = 0.3831459916571305
$ws.Range("M23").Value = 0.3583498223981678
$ws.Range("N23").Value = 1.725840481720525
$ws.Range("B24").Value = 1.282280990908134
$ws.Range("D24").Value = 0.006256287288714191
$ws.Range("E24").Value = 0.7857696892156838
$ws.Range("F24").Value = 1.1452018915152
$ws.Range("G24").Value = 1.079606747979511
$ws.Range("H24").Value = 0.8358158800868978
$ws.Range("L24").Value = 0.3362587112806921
$ws.Range("M24").Value = 0.3266723045427469
$ws.Range("N24").Value = 1.730359211313953
$ws.Range("B25").Value = 1.190483739365646
$ws.Range("D25").Value = 0.006249090091134946
$ws.Range("E25").Value = 0.7215454334010474
$ws.Range("F25").Value = 1.031660763430352
$ws.Range("G25").Value = 0.951526507928719
$ws.Range("H25").Value = 0.7863083565349314
$ws.Range("L25").Value = 0.2858987665469499
$ws.Range("M25").Value = 0.2929001545253342
$ws.Range("N25").Value = 1.737977092210045
